$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update header values (row 15) for columns Y..BR, then clear BS15 ---
$ws.Range("Y15").Value = 'ethnicity'
$ws.Range("Z15").Value = 'extrachrom_elements'
$ws.Range("AA15").Value = 'health_state'
$ws.Range("AB15").Value = 'host_age'
$ws.Range("AC15").Value = 'host_body_mass_index'
$ws.Range("AD15").Value = 'host_body_product'
$ws.Range("AE15").Value = 'host_body_temp'
$ws.Range("AF15").Value = 'host_diet'
$ws.Range("AG15").Value = 'host_disease'
$ws.Range("AH15").Value = 'host_family_relationship'
$ws.Range("AI15").Value = 'host_genotype'
$ws.Range("AJ15").Value = 'host_height'
$ws.Range("AK15").Value = 'host_last_meal'
$ws.Range("AL15").Value = 'host_occupation'
$ws.Range("AM15").Value = 'host_phenotype'
$ws.Range("AN15").Value = 'host_pulse'
$ws.Range("AO15").Value = 'host_sex'
$ws.Range("AP15").Value = 'host_subject_id'
$ws.Range("AQ15").Value = 'host_taxid'
$ws.Range("AR15").Value = 'host_tissue_sampled'
$ws.Range("AS15").Value = 'host_tot_mass'
$ws.Range("AT15").Value = 'ihmc_medication_code'
$ws.Range("AU15").Value = 'isolation_source'
$ws.Range("AV15").Value = 'locus_tag_prefix'
$ws.Range("AW15").Value = 'medic_hist_perform'
$ws.Range("AX15").Value = 'misc_param'
$ws.Range("AY15").Value = 'nose_mouth_teeth_throat_disord'
$ws.Range("AZ15").Value = 'organism_count'
$ws.Range("BA15").Value = 'oxy_stat_samp'
$ws.Range("BB15").Value = 'pathogenicity'
$ws.Range("BC15").Value = 'perturbation'
$ws.Range("BD15").Value = 'ref_biomaterial'
$ws.Range("BE15").Value = 'samp_collect_device'
$ws.Range("BF15").Value = 'samp_mat_process'
$ws.Range("BG15").Value = 'samp_salinity'
$ws.Range("BH15").Value = 'samp_size'
$ws.Range("BI15").Value = 'samp_store_dur'
$ws.Range("BJ15").Value = 'samp_store_loc'
$ws.Range("BK15").Value = 'samp_store_temp'
$ws.Range("BL15").Value = 'samp_vol_we_dna_ext'
$ws.Range("BM15").Value = 'source_material_id'
$ws.Range("BN15").Value = 'specimen_voucher'
$ws.Range("BO15").Value = 'subspecf_gen_lin'
$ws.Range("BP15").Value = 'temperature'
$ws.Range("BQ15").Value = 'time_last_toothbrush'
$ws.Range("BR15").Value = 'trophic_level'
$ws.Range("BS15").ClearContents()

# --- Step 2: Update comments for columns Y..BR (and remove the orphaned tail) ---
$ws.Range("Y15").Comment.Text('ethnicity of the subject')
$ws.Range("Z15").Comment.Text('Plasmids that have significance phenotypic consequence')
$ws.Range("AA15").Comment.Text('Health or disease status of sample at time of collection')
$ws.Range("AB15").Comment.Text('Age of host at the time of sampling')
$ws.Range("AC15").Comment.Text('body mass index of the host, calculated as weight/(height)squared')
$ws.Range("AD15").Comment.Text('substance produced by the host, e.g. stool, mucus, where the sample was obtained from')
$ws.Range("AE15").Comment.Text('core body temperature of the host when sample was collected')
$ws.Range("AF15").Comment.Text('type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types')
$ws.Range("AG15").Comment.Text('Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh')
$ws.Range("AH15").Comment.Delete()
$ws.Range("AJ15").AddComment('the height of subject')
$ws.Range("AK15").Comment.Text('content of last meal and time since feeding; can include multiple values')
$ws.Range("AL15").Comment.Text('most frequent job performed by subject')
$ws.Range("AM15").Comment.Delete()
$ws.Range("AN15").AddComment('resting pulse of the host, measured as beats per minute')
$ws.Range("AO15").Comment.Text('Gender or physical sex of the host')
$ws.Range("AP15").Comment.Text('a unique identifier by which each subject can be referred to, de-identified, e.g. #131')
$ws.Range("AQ15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606')
$ws.Range("AR15").Comment.Text('Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005')
$ws.Range("AS15").Comment.Text('total mass of the host at collection, the unit depends on host')
$ws.Range("AT15").Comment.Text('can include multiple medication codes')
$ws.Range("AU15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("AV15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$ws.Range("AW15").Comment.Text('whether full medical history was collected')
$ws.Range("AX15").Comment.Text('any other measurement performed or parameter collected, that is not listed here')
$ws.Range("AY15").Comment.Text('history of nose/mouth/teeth/throat disorders; can include multiple disorders')
$ws.Range("AZ15").Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$ws.Range("BA15").Comment.Text('oxygenation status of sample')
$ws.Range("BB15").Comment.Text('To what is the entity pathogenic')
$ws.Range("BC15").Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$ws.Range("BD15").Comment.Text('Primary publication or genome report in the form of pubmed ID, DOI or URL')
$ws.Range("BE15").Comment.Text('Method or device employed for collecting sample')
$ws.Range("BF15").Comment.Text('Processing applied to the sample during or after isolation')
$ws.Range("BG15").Comment.Text('salinity of sample, i.e. measure of total salt concentration')
$ws.Range("BH15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("BI15").Comment.Text('duration for which sample was stored')
$ws.Range("BJ15").Comment.Text('location at which sample was stored, usually name of a specific freezer/room')
$ws.Range("BK15").Comment.Text('temperature at which sample was stored, e.g. -80')
$ws.Range("BL15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("BM15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("BN15").Comment.Text('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier')
$ws.Range("BO15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$ws.Range("BP15").Comment.Text('temperature of the sample at time of sampling')
$ws.Range("BQ15").Comment.Text('specification of the time since last toothbrushing')
$ws.Range("BR15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)')
$ws.Range("BS15").Comment.Delete()
